$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add "url" in A1, shift "Endpoint" to B1 (already there),
#     add "Optional auth type" in C1 and "auth" in D1, all using the bold/border
#     header style already used by B1/C1. ---
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("D1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A1").Value = "url"
$ws.Range("B1").Value = "Endpoint"
$ws.Range("C1").Value = "Optional auth type"
$ws.Range("D1").Value = "auth"

# --- Data rows (2-4): the old numeric index column A (with header style and
#     numeric 0/1/2 values) is replaced by real URL strings with no special
#     style, and a new "auth" column D is introduced. ---
$ws.Range("A2:A4").ClearFormats() | Out-Null

$ws.Range("A2").Value = "http://127.0.0.1:8000/"
$ws.Range("B2").Value = "v1/example"
$ws.Range("C2").Value = "bearer"
$ws.Range("D2").Value = "Bearer xxxx"

$ws.Range("A3").Value = "http://127.0.0.1:8000/"
$ws.Range("B3").Value = "v1/example2"
$ws.Range("C3").Value = "basic"
$ws.Range("D3").Value = "cHJ1ZWJhMTIzMTM="

$ws.Range("A4").Value = "http://127.0.0.1:8000/123"
$ws.Range("B4").Value = "v1/example3"
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = ""
